# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Wed Apr  3 13:29:30 UTC 2024 with GitHub Actions"
# Updates price (col D) / 1h-volume-change (col E) text values, and fixes two
# row-ordering swaps (Maker/OKB at rows 38-39, ThetaToken/FirstDigitalUSD at rows 50-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.688.69'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '3.296.57'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '187.07'
$ws.Range('E5').Value = '  +5.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '552.69'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.577'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').Value = '3.287.69'
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.178'
$ws.Range('E10').Value = '  -1.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.579'
$ws.Range('E11').Value = '  +0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.52'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000264'
$ws.Range('E13').Value = '  +2.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.58'
$ws.Range('E14').Value = '  +1.90%  '
$ws.Range('D15').Value = '3.833.95'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '594.41'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '65.776.90'
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.80'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('D20').Value = '3.299.69'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.94'
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.893'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.42'
$ws.Range('E23').Value = '  +7.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.07'
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '100.41'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.92'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.04'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.72'
$ws.Range('E28').Value = '  +3.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.42'
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.62'
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.08'
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.66'
$ws.Range('E32').Value = '  +8.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.81'
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '571.05'
$ws.Range('E34').Value = '  +9.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '10.95'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.103'
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '56.85'
$ws.Range('E38').Value = '  +2.71%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.683.89'
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.45'
$ws.Range('E40').Value = '  +10.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '33.47'
$ws.Range('E41').Value = '  +6.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.22'
$ws.Range('E42').Value = '  -4.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.127'
$ws.Range('E43').Value = '  +2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.65'
$ws.Range('E44').Value = '  +1.54%  '
$ws.Range('D45').Value = '0.0₃0687'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.39'
$ws.Range('E46').Value = '  +6.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.335'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0414'
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.128'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.54'
$ws.Range('E51').Value = '  -0.18%  '
